$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Restructure columns: insert a new "Lake" column before the data, and a
# new "Description" column after the existing RMSE column. -----------------
$ws.Range("A1").EntireColumn.Insert()
$ws.Range("E1").EntireColumn.Insert()

# New column width for the Description column.
$ws.Columns.Item(5).ColumnWidth = 25.3

# --- Description column (entered first, so its strings land first in the
# shared-string table) ----------------------------------------------------
$ws.Range("B1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E1").Value2 = "Description"

$ws.Range("E2").Value2 = "No salt"
$ws.Range("E3").Value2 = "Constant salt value of 0.1"
$ws.Range("E4").Value2 = "Constant salt value of 01"
$ws.Range("E5").Value2 = "Constant salt value of 10"

# --- Lake column header + Mendota block -----------------------------------
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Range("A1").Value2 = "Lake"

$ws.Range("A2").Value2 = "Mendota"
$ws.Range("A3").Value2 = "Mendota"
$ws.Range("A4").Value2 = "Mendota"
$ws.Range("A5").Value2 = "Mendota"

# --- NSE results for the Mendota runs -------------------------------------
$ws.Range("C2").Value2 = 0.95792560000000004
$ws.Range("C3").Value2 = 0.95488709999999999
$ws.Range("C4").Value2 = 0.95688099999999998
$ws.Range("C5").Value2 = 0.8830692

# --- Monona block: new scenarios awaiting NSE/RMSE results ---------------
$ws.Range("A6").Value2 = "Monona"
$ws.Range("B6").Value2 = "Control"
$ws.Range("E6").Value2 = "No salt"

$ws.Range("A7").Value2 = "Monona"
$ws.Range("B7").Value2 = "A1"
$ws.Range("E7").Value2 = "Constant salt value of 0.1"

$ws.Range("A8").Value2 = "Monona"
$ws.Range("B8").Value2 = "A2"
$ws.Range("E8").Value2 = "Constant salt value of 01"

$ws.Range("A9").Value2 = "Monona"
$ws.Range("B9").Value2 = "A3"
$ws.Range("E9").Value2 = "Constant salt value of 10"

# --- Selection, matching the saved cursor position in the source file ---
$ws.Range("D8").Select()

Write-Output "done"
